# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" sheet (as the 2nd tab, right after "总计") with the
# latest quarterly fund-holdings data, and adds a corresponding summary row
# at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)      # "总计"
$wsQ3    = $wb.Worksheets.Item(2)      # "2021-Q3" - used as a formatting template

# --- 1. Create the new "2022-Q3" sheet by copying the "2021-Q3" sheet's ----
#        layout/formatting, placed right after "总计" -----------------------
$wsQ3.Copy($null, $wsTotal)
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q3"

# Header: "基金金额" -> "基金规模"
$wsNew.Range("D1").Value = "基金规模"

# Make sure the text-like columns (fund code / name / numbers-as-text) keep
# their text representation (preserve leading zeros, exact decimal text).
$wsNew.Range("B2:G3").NumberFormat = "@"

# Row 2 data
$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "005143"
$wsNew.Range("C2").Value = "中融沪港深大消费主题灵活配置混合C"
$wsNew.Range("D2").Value = "0.27"
$wsNew.Range("E2").Value = "90.10"
$wsNew.Range("F2").Value = "4.85"
$wsNew.Range("G2").Value = "0.0131"
$wsNew.Range("H2").Value = 7

# Row 3 is new - clone row 2's formatting down into row 3 first
$wsNew.Range("A2:H2").Copy()
$wsNew.Range("A3:H3").PasteSpecial(-4122)   # xlPasteFormats
$wsNew.Range("B3:G3").NumberFormat = "@"

$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "005142"
$wsNew.Range("C3").Value = "中融沪港深大消费主题灵活配置混合A"
$wsNew.Range("D3").Value = "0.13"
$wsNew.Range("E3").Value = "90.10"
$wsNew.Range("F3").Value = "4.85"
$wsNew.Range("G3").Value = "0.0063"
$wsNew.Range("H3").Value = 7

# --- 2. Add the "2022-Q3" summary row at the top of "总计"'s data ----------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)   # xlPasteFormats

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.02

# Re-assert the (now shifted) rows' index column + data exactly, so the
# sequential index in column A stays 0,1,2,3 top-to-bottom.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 3

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q2"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 5.09

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2021-Q1"
$wsTotal.Range("C5").Value = 1
$wsTotal.Range("D5").Value = 2.17

$wsNew.Select()
$wsNew.Range("A1").Select()
